$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the previous layout (old headers, spain-only row, country-only rows) ---
$ws.Hyperlinks.Delete()
$ws.Range("A1:F10").Clear()

# --- Headers (row 1) ---
$ws.Range("A1").Value = "Data Engineer"
$ws.Range("B1").Value = "Data Scientist"
$ws.Range("C1").Value = "Data Analyst"
$ws.Range("D1").Value = "Business Analyst"
$ws.Range("E1").Value = "Area"

# --- Spain row (row 2) ---
$ws.Range("A2").Value = "https://www.glassdoor.com/Job/spain-data-engineer-jobs-SRCH_IL.0,5_IN219_KO6,19.htm?clickSource=searchBox"
$ws.Range("B2").Value = "https://www.glassdoor.com/Job/spain-data-scientist-jobs-SRCH_IL.0,5_IN219_KO6,20.htm?clickSource=searchBox"
$ws.Range("C2").Value = "https://www.glassdoor.com/Job/spain-data-analyst-jobs-SRCH_IL.0,5_IN219_KO6,18.htm?clickSource=searchBox"
$ws.Range("D2").Value = "https://www.glassdoor.com/Job/spain-business-analyst-jobs-SRCH_IL.0,5_IN219_KO6,22.htm?clickSource=searchBox"
$ws.Range("E2").Value = "Espana"

# --- Barcelona row (row 3) ---
$ws.Range("A3").Value = "https://www.glassdoor.com/Job/barcelona-data-engineer-jobs-SRCH_IL.0,9_IC2547194_KO10,23.htm?clickSource=searchBox"
$ws.Range("B3").Value = "https://www.glassdoor.com/Job/barcelona-data-scientist-jobs-SRCH_IL.0,9_IC2547194_KO10,24.htm?clickSource=searchBox"
$ws.Range("C3").Value = "https://www.glassdoor.com/Job/barcelona-data-analyst-jobs-SRCH_IL.0,9_IC2547194_KO10,22.htm?clickSource=searchBox"
$ws.Range("D3").Value = "https://www.glassdoor.com/Job/barcelona-business-analyst-jobs-SRCH_IL.0,9_IC2547194_KO10,26.htm?clickSource=searchBox"
$ws.Range("E3").Value = "Barcelona"

# --- Madrid row (row 4) ---
$ws.Range("A4").Value = "https://www.glassdoor.com/Job/madrid-data-engineer-jobs-SRCH_IL.0,6_IC2664239_KO7,20.htm?clickSource=searchBox"
$ws.Range("B4").Value = "https://www.glassdoor.com/Job/madrid-data-scientist-jobs-SRCH_IL.0,6_IC2664239_KO7,21.htm?clickSource=searchBox"
$ws.Range("C4").Value = "https://www.glassdoor.com/Job/madrid-data-analyst-jobs-SRCH_IL.0,6_IC2664239_KO7,19.htm?clickSource=searchBox"
$ws.Range("D4").Value = "https://www.glassdoor.com/Job/madrid-business-analyst-jobs-SRCH_IL.0,6_IC2664239_KO7,23.htm?clickSource=searchBox"
$ws.Range("E4").Value = "Madrid"

# --- Hyperlinks, added in the same order the author created them ---
$ws.Hyperlinks.Add($ws.Range("A2"), $ws.Range("A2").Value2)
$ws.Hyperlinks.Add($ws.Range("B2"), $ws.Range("B2").Value2)
$ws.Hyperlinks.Add($ws.Range("C2"), $ws.Range("C2").Value2)
$ws.Hyperlinks.Add($ws.Range("D2"), $ws.Range("D2").Value2)
$ws.Hyperlinks.Add($ws.Range("B4"), $ws.Range("B4").Value2)
$ws.Hyperlinks.Add($ws.Range("B3"), $ws.Range("B3").Value2)
$ws.Hyperlinks.Add($ws.Range("C3"), $ws.Range("C3").Value2)
$ws.Hyperlinks.Add($ws.Range("C4"), $ws.Range("C4").Value2)
$ws.Hyperlinks.Add($ws.Range("D3"), $ws.Range("D3").Value2)

# Hyperlinks.Add stamps its own one-off style copy; collapse back onto the
# shared built-in "Hyperlink" cell style so linked cells share one style id.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"

# --- View / selection bookkeeping ---
$ws.Range("C1").Select()
$win = $excel.ActiveWindow
$win.Left = 2200
$win.Top = 2200
$win.Width = 14400
$win.Height = 7270
